$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 188 (shifts existing rows 188:305 down to 189:306,
# and grows the used range from A1:R305 to A1:R306).
$ws.Rows("188").Insert()

# Populate the newly inserted row 188 with the new weekly price record.
# Non-numeric / unchanged columns mirror the row immediately below it (old row 188,
# now row 189): Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría,
# Variedad, Calidad, Unidad de comercialización, Origen, Kg o Unidades, Clasificación.
$ws.Cells.Item(188, 1).Value = 7
$ws.Cells.Item(188, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(188, 3).Value = "Ñuble"
$ws.Cells.Item(188, 4).Value = 45161
$ws.Cells.Item(188, 5).Value = 16
$ws.Cells.Item(188, 6).Value = 100112040
$ws.Cells.Item(188, 7).Value = "Cilantro"
$ws.Cells.Item(188, 8).Value = "Sin especificar"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 300
$ws.Cells.Item(188, 11).Value = 1500
$ws.Cells.Item(188, 12).Value = 1500
$ws.Cells.Item(188, 13).Value = 1500
$ws.Cells.Item(188, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(188, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(188, 16).Value = 1500
$ws.Cells.Item(188, 17).Value = 1
$ws.Cells.Item(188, 18).Value = "Hortaliza"
